$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newest day's data after running the profit report on 2025-11-05.
# Existing rows go through row 79 (11/04/2025); the new data lands in row 80.
$row = 80
$dateText   = "11/05/2025"
$profitValue = 8875.940000000001

# The Date column stores its values as plain text (e.g. "11/04/2025" in A79),
# not as real Excel dates. If we simply did:
#   $ws.Range("A80").Value = $dateText
# Excel's automatic data-type detection would parse the date-looking string
# and store it as a date serial number instead of literal text. To avoid
# that, write it first as a text formula (so it is produced as a String
# result rather than being reinterpreted), then convert that formula cell
# to a plain value in place via copy/paste-special. This keeps the cell a
# normal text value with no special number format or style applied to it,
# matching the rest of the column.
$ws.Range("A$row").Formula = '="' + $dateText + '"'
$ws.Range("A$row").Copy()
$ws.Range("A$row").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Profit column is a normal numeric value.
$ws.Range("B$row").Value = $profitValue
